# Add the new worksheet "Spannungsmessung" after the existing "Strommessung" sheet
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Add([Type]::Missing, $ws1)
$ws2.Name = "Spannungsmessung"

# Row 3
$ws2.Range("A3").Value = "Eingangsspannung"
$ws2.Range("B3").Value = 120
$ws2.Range("C3").Value = "V"
$ws2.Range("D3").Value = 120
$ws2.Range("E3").Value = "V"

# Row 4
$ws2.Range("A4").Value = "Ausgangspannung"
$ws2.Range("B4").Value = 5
$ws2.Range("C4").Value = "V"

# Row 5
$ws2.Range("A5").Value = "R1"
$ws2.Range("B5").Formula = "=47000+4700"
$ws2.Range("C5").Value = "Ohm"
$ws2.Range("D5").Formula = "=47000+4700"
$ws2.Range("E5").Value = "Ohm"
$ws2.Range("F5").Formula = "=D3/(D5+D6)*D5"
$ws2.Range("G5").Value = "V"
$ws2.Range("H5").Formula = "=(F5*F5)/D5"

# Row 6
$ws2.Range("A6").Value = "R2"
$ws2.Range("B6").Formula = "=B5/(B3-B4)*B4"
$ws2.Range("C6").Value = "Ohm"
$ws2.Range("D6").Value = 2200
$ws2.Range("E6").Value = "Ohm"
$ws2.Range("F6").Formula = "=D3/(D5+D6)*D6"
$ws2.Range("G6").Value = "V"
$ws2.Range("H6").Formula = "=(F6*F6)/D6"

# Number formats - applied after all formulas so the propagated-format
# heuristic doesn't leak into dependent formula cells (H5/H6).
$ws2.Range("B5").NumberFormat = "#,##0"
$ws2.Range("B6").NumberFormat = "#,##0"
$ws2.Range("F5").NumberFormat = "#,##0.00"
$ws2.Range("F6").NumberFormat = "#,##0.00"

# Column widths to match target (col A bestFit width ~17.43, col B width ~11.86)
$ws2.Columns.Item(1).ColumnWidth = 17.42578125
$ws2.Columns.Item(2).ColumnWidth = 11.85546875

# Page setup to match target (A4 portrait, 2cm top/bottom margins)
$ws2.PageSetup.TopMargin = 56.692913399999995
$ws2.PageSetup.BottomMargin = 56.692913399999995
$ws2.PageSetup.PaperSize = 9
$ws2.PageSetup.Orientation = 1

# Select D5 like the target file, and make this the active/visible tab
$ws2.Range("D5").Select()

$wb.Save()
